# fix: alterar python version para 3.11.5
# Update absenteeism data rows 2-11 with new values as per upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=45473; B="Beatriz Souza";           C="TI";                 D="Viagem de negócios"; E=7; F=45085; G=6627.23}
    @{Row=3;  A=76086; B="Bruna Barbosa";            C="Marketing";          D="Viagem de negócios"; E=6; F=45103; G=5671.62}
    @{Row=4;  A=35171; B="Sr. Rodrigo Carvalho";     C="TI";                 D="Problemas pessoais"; E=4; F=45091; G=8180.57}
    @{Row=5;  A=36413; B="Juliana Souza";            C="Jurídico";           D="Viagem de negócios"; E=3; F=45093; G=11309.07}
    @{Row=6;  A=31466; B="Nina Lopes";                C="Recursos Humanos";  D="Consulta médica";    E=6; F=45094; G=8171.11}
    @{Row=7;  A=18287; B="João Miguel da Conceição"; C="Jurídico";           D="Doença";              E=5; F=45095; G=8191.68}
    @{Row=8;  A=71317; B="Srta. Vitória Nunes";      C="Engenharia";         D="Outros";              E=1; F=45079; G=5510.27}
    @{Row=9;  A=81764; B="Ana Beatriz da Rosa";      C="P&D";                D="Viagem de negócios"; E=8; F=45101; G=3496.78}
    @{Row=10; A=10726; B="Pedro Miguel Barbosa";     C="Jurídico";           D="Viagem de negócios"; E=8; F=45086; G=6562.34}
    @{Row=11; A=55286; B="Enzo Gabriel Sales";       C="Jurídico";           D="Problemas pessoais"; E=6; F=45100; G=7599.68}
)

foreach ($r in $data) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}
